# ScoreMatrix.xlsx edit:
#   - Cell B9 on Sheet1 held the note about the brand-match boost; the
#     percentage was corrected from 9% to 7%.
#   - The active selection/cursor was left sitting on B9 (the cell that was
#     just edited) instead of B10.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B9").Value = "7% boost to brand matches over 0, baseline model"

$ws.Range("B9").Select()

$wb.Save()
